$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their literal text representation
# (avoids Excel auto-converting numeric-looking strings and dropping trailing zeros)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '28.019.23'
$ws.Range('E2').Value = '  -2.22%  '
$ws.Range('D3').Value = '1.830.67'
$ws.Range('E3').Value = '  -1.20%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '324.26'
$ws.Range('E5').Value = '  -3.15%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D7').Value = '0.4640'
$ws.Range('E7').Value = '  -0.35%  '
$ws.Range('D8').Value = '0.3866'
$ws.Range('E8').Value = '  -1.21%  '
$ws.Range('D9').Value = '0.07835'
$ws.Range('E9').Value = '  -0.88%  '
$ws.Range('D10').Value = '0.9589'
$ws.Range('E10').Value = '  -2.65%  '
$ws.Range('D11').Value = '21.88'
$ws.Range('E11').Value = '  -2.08%  '
$ws.Range('D12').Value = '1.810.53'
$ws.Range('E12').Value = '  -2.21%  '
$ws.Range('D13').Value = '5.679'
$ws.Range('E13').Value = '  -3.08%  '
$ws.Range('D14').Value = '6.890'
$ws.Range('E14').Value = '  -1.78%  '
$ws.Range('D15').Value = '0.06866'
$ws.Range('E15').Value = '  -0.16%  '
$ws.Range('D16').Value = '88.19'
$ws.Range('E16').Value = '  +0.56%  '
$ws.Range('D17').Value = '1.003'
$ws.Range('E17').Value = '  +0.18%  '
$ws.Range('D18').Value = '0.000009911'
$ws.Range('E18').Value = '  -1.49%  '
$ws.Range('D19').Value = '16.63'
$ws.Range('E19').Value = '  -3.12%  '
$ws.Range('D20').Value = '1.000'
$ws.Range('D21').Value = '28.039.47'
$ws.Range('E21').Value = '  -2.21%  '
$ws.Range('D22').Value = '5.298'
$ws.Range('E22').Value = '  -1.98%  '
$ws.Range('D23').Value = '10.97'
$ws.Range('E23').Value = '  -3.41%  '
$ws.Range('D24').Value = '2.093'
$ws.Range('E24').Value = '  -2.27%  '
$ws.Range('D25').Value = '2.038.28'
$ws.Range('E25').Value = '  +0.81%  '
$ws.Range('D26').Value = '154.64'
$ws.Range('E26').Value = '  +0.80%  '
$ws.Range('D27').Value = '19.15'
$ws.Range('E27').Value = '  -1.85%  '
$ws.Range('D28').Value = '5.646'
$ws.Range('E28').Value = '  -6.75%  '
$ws.Range('D29').Value = '1.954'
$ws.Range('E29').Value = '  -3.92%  '
$ws.Range('D30').Value = '118.04'
$ws.Range('E30').Value = '  +0.26%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = '0.09236'
$ws.Range('E31').Value = '  -1.70%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = '0.9328'
$ws.Range('E32').Value = '  -4.61%  '
$ws.Range('D33').Value = '5.246'
$ws.Range('E33').Value = '  -2.37%  '
$ws.Range('D34').Value = '1.318'
$ws.Range('E34').Value = '  -2.45%  '
$ws.Range('E35').Value = '  -5.05%  '
$ws.Range('D36').Value = '0.05847'
$ws.Range('E36').Value = '  -5.09%  '
$ws.Range('D37').Value = '0.02124'
$ws.Range('E37').Value = '  -3.50%  '
$ws.Range('D38').Value = '1.145'
$ws.Range('E38').Value = '  -1.62%  '
$ws.Range('D39').Value = '7.734'
$ws.Range('E39').Value = '  +1.23%  '
$ws.Range('D40').Value = '0.5580'
$ws.Range('E40').Value = '  -2.70%  '
$ws.Range('D41').Value = '9.863'
$ws.Range('E41').Value = '  -3.39%  '
$ws.Range('D42').Value = '0.1758'
$ws.Range('E42').Value = '  -2.51%  '
$ws.Range('D43').Value = '0.07195'
$ws.Range('E43').Value = '  +0.60%  '
$ws.Range('D44').Value = '11.61'
$ws.Range('E44').Value = '  -1.23%  '
$ws.Range('D45').Value = '0.5258'
$ws.Range('E45').Value = '  -2.67%  '
$ws.Range('D46').Value = '1.145'
$ws.Range('E46').Value = '  -6.82%  '
$ws.Range('D47').Value = '2.096'
$ws.Range('E47').Value = '  -11.44%  '
$ws.Range('D48').Value = '1.820'
$ws.Range('E48').Value = '  -4.98%  '
$ws.Range('D49').Value = '112.20'
$ws.Range('E49').Value = '  -3.01%  '
$ws.Range('E50').Value = '  +0.15%  '
$ws.Range('B51').Value = 'MXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D51').Value = '2.328'
$ws.Range('E51').Value = '  +0.59%  '
